# Applies the "network_data.xlsx" update:
#  - Nodes sheet: rename P/Q headers, rescale column D (Q) by 5x
#  - Edges sheet: add a new D[m] diameter column, add a constant roughness column
#  - Params sheet: rename headers, change a few parameter values
#  - New "Daily demand" sheet with hourly demand factors

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 0. Create the new "Daily demand" sheet and position it after "Params" FIRST
#    (sheet handles captured before a Move() go stale, so do this before
#    grabbing references to the other sheets we need to edit).
# ---------------------------------------------------------------------------
$daily = $wb.Worksheets.Add()
$daily.Name = "Daily demand"
$daily.Move($null, $wb.Worksheets.Item("Params"))

# Re-fetch every sheet reference fresh, by name, now that ordering is final.
$nodes = $wb.Worksheets.Item("Nodes")
$edges = $wb.Worksheets.Item("Edges")
$params = $wb.Worksheets.Item("Params")
$daily = $wb.Worksheets.Item("Daily demand")

# ---------------------------------------------------------------------------
# 1. Nodes sheet
# ---------------------------------------------------------------------------
$nodes.Range("C1").Value = "P[Pa]"
$nodes.Range("D1").Value = "Q[m^3/s]"

for ($r = 2; $r -le 15; $r++) {
    $cell = $nodes.Cells.Item($r, 4)
    $cell.Value = ($cell.Value2 * 5)
}

# ---------------------------------------------------------------------------
# 2. Edges sheet
# ---------------------------------------------------------------------------
$edges.Range("E1").Value = "D[m]"

$edges.Range("D2").Value = 2
for ($r = 3; $r -le 19; $r++) {
    $edges.Cells.Item($r, 4).Value = 1.3
}
for ($r = 2; $r -le 19; $r++) {
    $edges.Cells.Item($r, 5).Value = 0.3
}

# ---------------------------------------------------------------------------
# 3. Params sheet
# ---------------------------------------------------------------------------
$params.Range("A1").Value = "D [m]"
$params.Range("E1").Value = "nodes_with_boundary_conditions"

$params.Range("A2").Value = 0.7
$params.Range("B2").Value = 0.05
$params.Range("C2").Value = 0.00001357

$params.Range("A5").Value = "note: Do not initialize D. It is automatically calculated in the code."

# ---------------------------------------------------------------------------
# 4. "Daily demand" sheet contents
# ---------------------------------------------------------------------------
$daily.Range("A1").Value = "Q[m^3/s]"
$dailyValues = @(0.972222, 1.041667, 1.125, 1.388889, 1.666667, 1.805556, 1.597222, 1.458333, 1.597222, 1.666667, 1.458333, 0.972222)
for ($i = 0; $i -lt $dailyValues.Length; $i++) {
    $daily.Cells.Item($i + 2, 1).Value = $dailyValues[$i]
}

$daily.Activate()
$daily.Select()
